# Update summary-table statistic values and the footnote participant counts
# ("update to light glasses" - refreshed analysis results).
$d = $word.ActiveDocument

$replacements = @(
    ,@("<strong>1</strong> ", "<strong>20</strong> ")
    ,@("<strong>8</strong> (8 - 8)", "<strong>153</strong> (1 - 8)")
    ,@("<strong>6</strong> (6 - 6)", "<strong>115</strong> (6 - 7)")
    ,@("<strong>14.0%</strong> (14.0% - 14.0%)", "<strong>13.0%</strong> (7.0% - 64.0%)")
    ,@("<strong>13h 32s</strong> (12h 46m - 13h 14m)", "<strong>16h 31m</strong> (11h 31m - 18h 13m)")
    ,@("<strong>37</strong> ±90 (0 - 221)", "<strong>12,325</strong> ±21,523 (0 - 144,958)")
    ,@("<strong>2m 40s</strong> ±6m 31s (0s - 16m)", "<strong>3h 42m</strong> ±3h 9m (0s - 15h 29m)")
    ,@("<strong>13m 15s</strong> ±32m 27s (0s - 1h 19m)", "<strong>2h 27m</strong> ±1h 53m (0s - 10h 18m)")
    ,@("<strong>23h 38m</strong> ±52m 11s (21h 52m - 1d)", "<strong>10h 19m</strong> ±3h 54m (4h 51m - 1d)")
    ,@("<strong>2m 33s</strong> ±6m 15s (0s - 15m 20s)", "<strong>56m 57s</strong> ±1h 13m (0s - 6h 44m)")
    ,@("<strong>0s</strong> ±0s (0s - 0s)", "<strong>1h 37m</strong> ±2h 9m (0s - 14h 4m)")
    ,@("<strong>07:51</strong> (07:51 - 07:51)", "<strong>08:50</strong> ±02:25 (00:43 - 16:49)")
    ,@("<strong>09:17</strong> (09:17 - 09:17)", "<strong>13:51</strong> ±01:52 (09:17 - 19:35)")
    ,@("<strong>13:07</strong> (13:07 - 13:07)", "<strong>19:27</strong> ±02:17 (10:40 - 22:56)")
    ,@("<strong>05:52</strong> ±02:07 (04:59 - 10:12)", "<strong>13:24</strong> ±02:34 (04:59 - 18:59)")
    ,@("<strong>02:29</strong> ±00:00 (02:29 - 02:29)", "<strong>03:01</strong> ±02:55 (00:18 - 23:32)")
    ,@("<strong>0.0</strong> ±0.1 (0.0 - 0.2)", "<strong>320.2</strong> ±1,032.6 (0.0 - 9,111.6)")
    ,@("<strong>0.0</strong> ±0.0 (0.0 - 0.0)", "<strong>0.0</strong> ±0.0 (0.0 - 0.1)")
    ,@("<strong>0.158</strong> (0.158 - 0.158)", "<strong>0.264</strong> ±0.078 (0.158 - 0.414)")
    ,@("<strong>1.787</strong> (1.787 - 1.787)", "<strong>1.335</strong> ±0.414 (0.415 - 1.814)")
    ,@("Metrics are calculated on a by-participant-day basis (n=6) with the exception of IV and IS, which are calculated on a by-participant basis (n=1).", "Metrics are calculated on a by-participant-day basis (n=115) with the exception of IV and IS, which are calculated on a by-participant basis (n=20).")
)

foreach ($pair in $replacements) {
    $findText = $pair[0]
    $replaceText = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false,
                                 $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $findText"
    }
}
